$d = $word.ActiveDocument

# Locate the paragraph that still holds the old "still trying to get in
# contact..." placeholder text (it is the final paragraph of the
# "Acceptance Tests" section).
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*still trying to get in contact*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    Write-Output "Target paragraph not found!"
} else {
    $target = $d.Paragraphs.Item($targetIndex)

    # Insert a brand-new, completely bare paragraph (<w:p/>) right before it.
    $insPoint = $d.Range($target.Range.Start, $target.Range.Start)
    $insPoint.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

    # The original paragraph has now shifted down by one position.
    $target = $d.Paragraphs.Item($targetIndex + 1)

    # Replace its contents: drop the old sentence, keep the bold space run,
    # and append the new "Needs to be implemented." run (not bold). Also
    # drop the bold run-properties that used to live on the paragraph mark.
    $newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>' +
              '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
              '<w:r><w:t>Needs to be implemented.</w:t></w:r>' +
              '</w:p>'
    $target.Range.InsertXML($newXml)
}
